$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price/Volume columns store numeric-looking values and percentages as
# literal text (e.g. "328.80", "1.07%") rather than numbers. Force each cell
# being updated to Text format first so Excel does not reinterpret the new
# value as a number, preserving the original text data type.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.94%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.30%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.489"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.02%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08006"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.74%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.984"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.73%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.19%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9494"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.89%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1124"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.48%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1889"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.25%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "10.70"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "27.11%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1005"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.17%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04799"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12.79%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1060"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.33%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001273"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.78%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04081"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.02%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005983"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.59%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.17%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.374"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.29%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3472"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.74%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.71%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.88%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001268"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.57%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004332"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.95%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.76%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003746"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.02%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02581"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.92%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05657"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.51%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007559"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.31%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1397"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.17%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007402"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.65%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002016"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.32%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008645"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.80%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007111"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.36%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.13%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003532"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "55.71%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003647"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.80%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.13%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.13%"
